$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "70.084.43"
$ws.Range("E2").Value = "  -1.29%  "
$ws.Range("D3").Value = "3.573.08"
$ws.Range("E3").Value = "  -2.00%  "
$ws.Range("E4").Value = "  +0.17%  "
$ws.Range("D5").Value = "'576.82"
$ws.Range("E5").Value = "  -2.90%  "
$ws.Range("D6").Value = "'186.60"
$ws.Range("E6").Value = "  -4.38%  "
$ws.Range("D7").Value = "3.569.16"
$ws.Range("E7").Value = "  -1.95%  "
$ws.Range("D8").Value = "'0.621"
$ws.Range("E8").Value = "  -3.62%  "
$ws.Range("E9").Value = "  +0.12%  "
$ws.Range("E10").Value = "  +1.29%  "
$ws.Range("D11").Value = "'0.653"
$ws.Range("E11").Value = "  -3.48%  "
$ws.Range("D12").Value = "'54.44"
$ws.Range("E12").Value = "  -6.03%  "
$ws.Range("D13").Value = "'0.0000304"
$ws.Range("E13").Value = "  -2.09%  "
$ws.Range("D14").Value = "'9.56"
$ws.Range("E14").Value = "  -3.77%  "
$ws.Range("D15").Value = "4.157.33"
$ws.Range("E15").Value = "  -1.75%  "
$ws.Range("D16").Value = "'19.66"
$ws.Range("E16").Value = "  -4.40%  "
$ws.Range("D17").Value = "3.583.51"
$ws.Range("E17").Value = "  -1.68%  "
$ws.Range("D18").Value = "70.104.16"
$ws.Range("E18").Value = "  -1.19%  "
$ws.Range("D19").Value = "'12.47"
$ws.Range("E19").Value = "  -2.58%  "
$ws.Range("E20").Value = "  -1.40%  "
$ws.Range("E21").Value = "  -2.26%  "
$ws.Range("D22").Value = "'494.09"
$ws.Range("E22").Value = "  -0.05%  "
$ws.Range("D23").Value = "'19.37"
$ws.Range("E23").Value = "  +3.14%  "
$ws.Range("D24").Value = "'5.03"
$ws.Range("E24").Value = "  -3.56%  "
$ws.Range("D25").Value = "'96.64"
$ws.Range("E25").Value = "  +6.04%  "
$ws.Range("D26").Value = "'4.38"
$ws.Range("E26").Value = "  -3.53%  "
$ws.Range("D27").Value = "'11.59"
$ws.Range("E27").Value = "  +0.83%  "
$ws.Range("D28").Value = "'2.98"
$ws.Range("E28").Value = "  -5.98%  "
$ws.Range("D29").Value = "'9.34"
$ws.Range("E29").Value = "  -2.81%  "
$ws.Range("D30").Value = "'7.74"
$ws.Range("E30").Value = "  -2.85%  "
$ws.Range("D31").Value = "'31.61"
$ws.Range("E31").Value = "  -3.76%  "
$ws.Range("D32").Value = "'12.78"
$ws.Range("E32").Value = "  +3.79%  "
$ws.Range("D33").Value = "'65.58"
$ws.Range("E33").Value = "  -3.56%  "
$ws.Range("E34").Value = "  -5.16%  "
$ws.Range("D35").Value = "'573.78"
$ws.Range("E35").Value = "  -6.83%  "
$ws.Range("D36").Value = "'3.27"
$ws.Range("E36").Value = "  +11.30%  "
$ws.Range("D37").Value = "'38.89"
$ws.Range("E37").Value = "  -3.72%  "
$ws.Range("D38").Value = "'0.409"
$ws.Range("E38").Value = "  -0.68%  "
$ws.Range("E39").Value = "  -0.06%  "
$ws.Range("D40").Value = "0.0₃0792"
$ws.Range("E40").Value = "  -6.72%  "
$ws.Range("B41").Value = "ApeXProtocol"
$ws.Range("C41").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D41").Value = "'3.75"
$ws.Range("E41").Value = "  +9.91%  "
$ws.Range("B42").Value = "dogwifhat"
$ws.Range("C42").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D42").Value = "'3.18"
$ws.Range("E42").Value = "  -2.07%  "
$ws.Range("B43").Value = "Stacks"
$ws.Range("C43").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D43").Value = "'3.45"
$ws.Range("E43").Value = "  -3.54%  "
$ws.Range("E44").Value = "  -8.92%  "
$ws.Range("D45").Value = "'3.04"
$ws.Range("E45").Value = "  -3.78%  "
$ws.Range("D46").Value = "'0.0452"
$ws.Range("E46").Value = "  -1.91%  "
$ws.Range("D47").Value = "3.237.36"
$ws.Range("E47").Value = "  -2.79%  "
$ws.Range("D48").Value = "'9.51"
$ws.Range("E48").Value = "  -1.92%  "
$ws.Range("E49").Value = "  -2.76%  "
$ws.Range("D50").Value = "'1.48"
$ws.Range("E50").Value = "  +22.88%  "
$ws.Range("E51").Value = "  +0.13%  "
